$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value for column F (dSF), per repulled data
$updates = @{
    "F2"  = 2
    "F9"  = 2
    "F10" = -2
    "F15" = -3
    "F28" = -1
    "F29" = -4
    "F36" = 0
    "F38" = 2
    "F39" = 1
    "F41" = 3
    "F49" = -2
    "F52" = -6
    "F55" = -3
    "F57" = 0
    "F60" = 1
    "F65" = -1
    "F76" = 0
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
